$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the footer ("Nota"/"Fuente") rows - this pushes the
# existing row 21 ("Nota: El concepto...") down to row 22 and row 22
# ("Fuente: SICT...") down to row 23.
$ws.Rows("21:21").Insert()

# The new row 21 becomes a standalone "Actualización: mayo 2024" label,
# replacing the old note that used to live in N21/N22 ("Ultima actualización:
# mayo 2024" / "Dirección General de Planeación"). Give it the same look as
# the other footer labels (B22) and a numeric format across the rest of the
# row so it matches the table's body formatting.
$ws.Range("B22").Copy()
$ws.Range("B21").PasteSpecial(-4122)
$ws.Range("B21").Value2 = "Actualización: mayo 2024"
$ws.Range("C21:N21").NumberFormat = "#,##0.0"

# Inserting the row duplicated the old note text into N22/N23 (it rode along
# with the rows it was part of) - clear it out since that note is replaced by
# the new B21 label above, leaving the existing cell formatting untouched.
$ws.Range("N22").ClearContents()
$ws.Range("N23").ClearContents()
